# Generate Report for Handoff
# Update "Latest HO Xliff Generate Date" on Overview and Priority / Latest
# Handoff Datetime on the zh-cn and de-de language sheets for the rows
# whose handoff batch was just (re)generated.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 11, 12, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-09-07 10:29:39"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-07 10:29:34"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-07 10:29:39"
}
